$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# New row (row 5) on the "Library_Formula" sheet: add the getExpressionQuery
# BR detail line (IND_97 - add error column / BR details).
$ws.Range("B5").Value = "LIB_EWS"
$ws.Range("C5").Value = "getExpressionQuery"
$ws.Range("E5").Value = "String"
$ws.Range("F5").Value = "String,String"

$fmt = $ws.Range("E5:F5")
$fmt.Font.Name = "Trebuchet MS"
$fmt.Font.Size = 10
$fmt.Font.Color = 0

$ws.Range("F12").Select()
